$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 2: "Version number is written as the format: major.minor.build.[revision]"
# Add spellStart/spellEnd proofErr markers around "major.minor.build", change ".[" -> ". [",
# and drop the old gramStart/gramEnd proofErr markers.
$p2xml = "<w:p $wNs>" + `
  "<w:pPr>" + `
    "<w:spacing w:line='240' w:lineRule='auto'/>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
  "</w:pPr>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t xml:space='preserve'>Version number is written as the format: </w:t>" + `
  "</w:r>" + `
  "<w:proofErr w:type='spellStart'/>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:b/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t>major.minor.build</w:t>" + `
  "</w:r>" + `
  "<w:proofErr w:type='spellEnd'/>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:b/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t>. [</w:t>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:b/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t>revision]</w:t>" + `
  "</w:r>" + `
  "</w:p>"
$d.Paragraphs(2).Range.InsertXML($p2xml)

# --- Paragraph 3: "major: ... revision: Edit. Marked revision of the code." -- drop trailing bookmark.
$p3xml = "<w:p $wNs>" + `
  "<w:pPr>" + `
    "<w:spacing w:line='240' w:lineRule='auto'/>" + `
    "<w:ind w:left='720'/>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
  "</w:pPr>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:t>major: the major version, 1 digit.</w:t>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:br/>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:t>minor: the spinoff, 1 digit.</w:t>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:br/>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:t>build: version structure. Marked differences in the minor version, 1 or 2 digits.</w:t>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:br/>" + `
  "</w:r>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:i/>" + `
    "</w:rPr>" + `
    "<w:t>revision: Edit. Marked revision of the code.</w:t>" + `
  "</w:r>" + `
  "</w:p>"
$d.Paragraphs(3).Range.InsertXML($p3xml)

# --- Paragraphs 5, 7, 9: add justify-both alignment.
$d.Paragraphs(5).Range.ParagraphFormat.Alignment = 3
$d.Paragraphs(7).Range.ParagraphFormat.Alignment = 3
$d.Paragraphs(9).Range.ParagraphFormat.Alignment = 3

# --- Paragraph 11: add justify-both, split text and re-insert the _GoBack bookmark mid-run.
$p11xml = "<w:p $wNs>" + `
  "<w:pPr>" + `
    "<w:spacing w:line='240' w:lineRule='auto'/>" + `
    "<w:jc w:val='both'/>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
  "</w:pPr>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t xml:space='preserve'>- There is any change in source code that is not necessary to </w:t>" + `
  "</w:r>" + `
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" + `
  "<w:bookmarkEnd w:id='0'/>" + `
  "<w:r>" + `
    "<w:rPr>" + `
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
      "<w:sz w:val='24'/>" + `
    "</w:rPr>" + `
    "<w:t>change the name version. This index is a revision (revisions) of the source code, it marks the revision number of the source code (equivalent to GitHub commit)</w:t>" + `
  "</w:r>" + `
  "</w:p>"
$d.Paragraphs(11).Range.InsertXML($p11xml)
